$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 227.85715
$ws.Range("I6").Value = 227.85715
$ws.Range("K6").Value = 683.5714499999999
$ws.Range("M6").Value = -571.5714499999999
$ws.Range("H15").Value = 1455.0107
$ws.Range("I15").Value = 1455.0107
$ws.Range("K15").Value = 4365.0321
$ws.Range("M15").Value = -4196.0321
$ws.Range("H17").Value = 3449884.8
$ws.Range("J17").Value = 3847823
$ws.Range("L17").Value = 11543469
$ws.Range("N17").Value = -11543805
$ws.Range("H116").Value = 4488.727
$ws.Range("J116").Value = 4488.727
$ws.Range("L116").Value = 4488.727
$ws.Range("N116").Value = -11372.727
$ws.Range("H129").Value = 501386.84
$ws.Range("J129").Value = 589808.6
$ws.Range("L129").Value = 1769425.8
$ws.Range("N129").Value = -1779425.8
$ws.Range("H132").Value = 1871.6383
$ws.Range("I132").Value = 1988.341
$ws.Range("J132").Value = 160
$ws.Range("K132").Value = 5965.022999999999
$ws.Range("L132").Value = 480
$ws.Range("M132").Value = -3435.022999999999
$ws.Range("N132").Value = -5540
$ws.Range("H135").Value = 21745890
$ws.Range("I135").Value = 665.2222
$ws.Range("J135").Value = 100028700
$ws.Range("K135").Value = 5986.999800000001
$ws.Range("L135").Value = 900258300
$ws.Range("M135").Value = -3451.999800000001
$ws.Range("N135").Value = -900263370
$ws.Range("H137").Value = 1854.9166
$ws.Range("I137").Value = 1467.3077
$ws.Range("K137").Value = 4401.9231
$ws.Range("M137").Value = -1851.9231
$ws.Range("H138").Value = 2373.5208
$ws.Range("J138").Value = 2416
$ws.Range("L138").Value = 7248
$ws.Range("N138").Value = -17528
$ws.Range("H141").Value = 1185.9512
$ws.Range("I141").Value = 925.9722
$ws.Range("J141").Value = 3057.8
$ws.Range("K141").Value = 2777.9166
$ws.Range("L141").Value = 9173.400000000001
$ws.Range("M141").Value = 2402.0834
$ws.Range("N141").Value = -19533.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6164.4897
$ws.Range("I32").Value = 6443
$ws.Range("J32").Value = 4926.6665
$ws.Range("K32").Value = 6443
$ws.Range("L32").Value = 4926.6665
$ws.Range("M32").Value = -6156
$ws.Range("N32").Value = -5500.6665
$ws.Range("H41").Value = 3685.3333
$ws.Range("I41").Value = 3685.3333
$ws.Range("K41").Value = 3685.3333
$ws.Range("M41").Value = -3271.3333
$ws.Range("H61").Value = 1625.1428
$ws.Range("I61").Value = 1432.25
$ws.Range("J61").Value = 2782.5
$ws.Range("K61").Value = 1432.25
$ws.Range("L61").Value = 2782.5
$ws.Range("M61").Value = -1220.25
$ws.Range("N61").Value = -3206.5
$ws.Range("H74").Value = 71429330
$ws.Range("I74").Value = 76923850
$ws.Range("J74").Value = 600
$ws.Range("K74").Value = 76923850
$ws.Range("L74").Value = 600
$ws.Range("M74").Value = -76922976
$ws.Range("N74").Value = -2348
$ws.Range("H77").Value = 71429330
$ws.Range("I77").Value = 76923850
$ws.Range("J77").Value = 600
$ws.Range("K77").Value = 384619250
$ws.Range("L77").Value = 3000
$ws.Range("M77").Value = -384614882
$ws.Range("N77").Value = -11736
$ws.Range("H110").Value = 629.7143
$ws.Range("I110").Value = 634.6667
$ws.Range("J110").Value = 600
$ws.Range("K110").Value = 634.6667
$ws.Range("L110").Value = 600
$ws.Range("M110").Value = 1410.3333
$ws.Range("N110").Value = -4690
$ws.Range("H122").Value = 2034.1666
$ws.Range("I122").Value = 1785.25
$ws.Range("J122").Value = 2532
$ws.Range("K122").Value = 5355.75
$ws.Range("L122").Value = 7596
$ws.Range("M122").Value = -2905.75
$ws.Range("N122").Value = -12496
$ws.Range("H132").Value = 29530.445
$ws.Range("I132").Value = 1460.8605
$ws.Range("J132").Value = 139257
$ws.Range("K132").Value = 4382.5815
$ws.Range("L132").Value = 417771
$ws.Range("M132").Value = -1852.5815
$ws.Range("N132").Value = -422831
$ws.Range("H136").Value = 1625.1428
$ws.Range("I136").Value = 1432.25
$ws.Range("J136").Value = 2782.5
$ws.Range("K136").Value = 4296.75
$ws.Range("L136").Value = 8347.5
$ws.Range("M136").Value = -1746.75
$ws.Range("N136").Value = -13447.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7909.857
$ws.Range("I134").Value = 9111.272000000001
$ws.Range("K134").Value = 27333.816
$ws.Range("M134").Value = -24798.816

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11635.881
$ws.Range("I31").Value = 16463.154
$ws.Range("J31").Value = 3791.5625
$ws.Range("K31").Value = 16463.154
$ws.Range("L31").Value = 3791.5625
$ws.Range("M31").Value = -16168.154
$ws.Range("N31").Value = -4381.5625
$ws.Range("H34").Value = 11635.881
$ws.Range("I34").Value = 16463.154
$ws.Range("J34").Value = 3791.5625
$ws.Range("K34").Value = 16463.154
$ws.Range("L34").Value = 3791.5625
$ws.Range("M34").Value = -16261.154
$ws.Range("N34").Value = -4195.5625
$ws.Range("H107").Value = 1241.1
$ws.Range("I107").Value = 915.3333
$ws.Range("J107").Value = 1507.6364
$ws.Range("K107").Value = 915.3333
$ws.Range("L107").Value = 1507.6364
$ws.Range("M107").Value = 1004.6667
$ws.Range("N107").Value = -5347.6364
$ws.Range("H132").Value = 18195.25
$ws.Range("I132").Value = 21777.68
$ws.Range("J132").Value = 5400.857
$ws.Range("K132").Value = 65333.04
$ws.Range("L132").Value = 16202.571
$ws.Range("M132").Value = -62803.04
$ws.Range("N132").Value = -21262.571
$ws.Range("H134").Value = 651.1539
$ws.Range("I134").Value = 548.36664
$ws.Range("J134").Value = 993.7778
$ws.Range("K134").Value = 1645.09992
$ws.Range("L134").Value = 2981.3334
$ws.Range("M134").Value = 889.9000800000001
$ws.Range("N134").Value = -8051.3334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2777.625
$ws.Range("I3").Value = 1474.2727
$ws.Range("J3").Value = 5645
$ws.Range("K3").Value = 4422.8181
$ws.Range("L3").Value = 16935
$ws.Range("M3").Value = -4310.8181
$ws.Range("N3").Value = -17159
$ws.Range("H49").Value = 1750
$ws.Range("J49").Value = 1750
$ws.Range("L49").Value = 5250
$ws.Range("N49").Value = -5562
$ws.Range("H131").Value = 125816.36
$ws.Range("I131").Value = 486
$ws.Range("J131").Value = 134171.72
$ws.Range("K131").Value = 1458
$ws.Range("L131").Value = 402515.16
$ws.Range("M131").Value = 3582
$ws.Range("N131").Value = -412595.16
$ws.Range("H133").Value = 3763
$ws.Range("I133").Value = 1947.1428
$ws.Range("J133").Value = 8000
$ws.Range("K133").Value = 5841.428400000001
$ws.Range("L133").Value = 24000
$ws.Range("M133").Value = -781.4284000000007
$ws.Range("N133").Value = -34120
$ws.Range("H134").Value = 2175
$ws.Range("I134").Value = 1065.84
$ws.Range("J134").Value = 5641.125
$ws.Range("K134").Value = 3197.52
$ws.Range("L134").Value = 16923.375
$ws.Range("M134").Value = 1872.48
$ws.Range("N134").Value = -27063.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9461.388999999999
$ws.Range("I70").Value = 12447
$ws.Range("J70").Value = 4769.7144
$ws.Range("K70").Value = 12447
$ws.Range("L70").Value = 4769.7144
$ws.Range("M70").Value = -12177
$ws.Range("N70").Value = -5309.7144
$ws.Range("H73").Value = 9461.388999999999
$ws.Range("I73").Value = 12447
$ws.Range("J73").Value = 4769.7144
$ws.Range("K73").Value = 12447
$ws.Range("L73").Value = 4769.7144
$ws.Range("M73").Value = -11511
$ws.Range("N73").Value = -6641.7144
$ws.Range("H126").Value = 4128.485
$ws.Range("I126").Value = 3349.8096
$ws.Range("K126").Value = 10049.4288
$ws.Range("M126").Value = -7579.4288
$ws.Range("H132").Value = 18207.787
$ws.Range("I132").Value = 3211.6667
$ws.Range("J132").Value = 58197.445
$ws.Range("K132").Value = 9635.000100000001
$ws.Range("L132").Value = 174592.335
$ws.Range("M132").Value = -7105.000100000001
$ws.Range("N132").Value = -179652.335

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3682.7
$ws.Range("I40").Value = 3003.6667
$ws.Range("K40").Value = 3003.6667
$ws.Range("M40").Value = -2867.6667
$ws.Range("H132").Value = 1597.1562
$ws.Range("I132").Value = 1307.6296
$ws.Range("J132").Value = 3160.6
$ws.Range("K132").Value = 3922.8888
$ws.Range("L132").Value = 9481.799999999999
$ws.Range("M132").Value = -1392.8888
$ws.Range("N132").Value = -14541.8
$ws.Range("H136").Value = 35257
$ws.Range("I136").Value = 51374.5
$ws.Range("J136").Value = 3022
$ws.Range("K136").Value = 154123.5
$ws.Range("L136").Value = 9066
$ws.Range("M136").Value = -151573.5
$ws.Range("N136").Value = -14166

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 111113176
$ws.Range("I81").Value = 2183.5
$ws.Range("K81").Value = 4367
$ws.Range("M81").Value = -3306
$ws.Range("H84").Value = 111113176
$ws.Range("I84").Value = 2183.5
$ws.Range("K84").Value = 21835
$ws.Range("M84").Value = -16531
$ws.Range("H132").Value = 1228.2084
$ws.Range("I132").Value = 793.6316
$ws.Range("J132").Value = 2879.6
$ws.Range("K132").Value = 2380.8948
$ws.Range("L132").Value = 8638.799999999999
$ws.Range("M132").Value = 149.1052
$ws.Range("N132").Value = -13698.8
$ws.Range("H136").Value = 34484450
$ws.Range("I136").Value = 41668240
$ws.Range("J136").Value = 2260
$ws.Range("K136").Value = 125004720
$ws.Range("L136").Value = 6780
$ws.Range("M136").Value = -125002170
$ws.Range("N136").Value = -11880

Write-Host "Applied all Typhon_Profits updates"